$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the title row (old row 1); everything below shifts up by one.
$ws.Rows("1:1").Delete()

# Select row 2 (the header row), matching the saved selection state.
$ws.Range("A2:XFD2").Select()
